# "Generate Report for Handoff"
#
# Re-stamps the localization-status report with a new handoff commit
# (new source-file UUID "6a1047e3-72e5-4c12-ab6e-21fda646bcc0" and a new
# handoff bundle hash "aae9bb3b53f0d990fb5a784db8f4485271b272f2"), and
# bumps the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps to
# the new handoff run's time.
#
# NOTE: the hyperlink *targets* (the github blob URLs, captured at the
# time of the previous handoff) are intentionally left as-is -- only the
# visible/displayed text of each hyperlink is refreshed to match the new
# file name, mirroring the new cell text.

$wb = $excel.ActiveWorkbook

$oldUuid = "6352c826-2c1a-41f7-b6ed-e756ed4d7fea"
$newUuid = "6a1047e3-72e5-4c12-ab6e-21fda646bcc0"
$oldHash = "0f771765c78f2d806d97d84afdc7a7ed7937a3a3"
$newHash = "aae9bb3b53f0d990fb5a784db8f4485271b272f2"

# ---------------------------------------------------------------------
# Sheet "Overview": A2 = "<uuid>.md" (hyperlinked), D2 = handoff date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewA2Target = $wsOverview.Range("A2").Hyperlinks.Item(1).Address
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $overviewA2Target, [Type]::Missing, [Type]::Missing, "$newUuid.md") | Out-Null

$wsOverview.Range("D2").Value = "2016-55-13 02:55:52"

# ---------------------------------------------------------------------
# Sheet "zh-cn": A2 = "<uuid>.md", B2 = ".md" (unchanged), D2 = handoff
# xlf file name, E2 = handoff datetime
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhA2Target = $wsZh.Range("A2").Hyperlinks.Item(1).Address
$zhB2Target = $wsZh.Range("B2").Hyperlinks.Item(1).Address
$zhB2Display = "{0}" -f ".md"
$zhD2Target = $wsZh.Range("D2").Hyperlinks.Item(1).Address

# One Delete() on any cell's Hyperlinks clears the whole sheet's
# collection in this host, so grab every target first, then rebuild all
# three links (in their original rId order) with the refreshed display
# text for A2/D2 and the untouched display text for B2.
$wsZh.Range("A2").Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Target, [Type]::Missing, [Type]::Missing, "$newUuid.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $zhB2Target, [Type]::Missing, [Type]::Missing, $zhB2Display) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhD2Target, [Type]::Missing, [Type]::Missing, "$newUuid.$newHash.zh-cn.xlf") | Out-Null

$wsZh.Range("E2").Value = "2016-03-13 02:55:48"

# ---------------------------------------------------------------------
# Sheet "de-de": A2 = "<uuid>.md", B2 = ".md" (unchanged), D2 = handoff
# xlf file name, E2 = handoff datetime
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deA2Target = $wsDe.Range("A2").Hyperlinks.Item(1).Address
$deB2Target = $wsDe.Range("B2").Hyperlinks.Item(1).Address
$deB2Display = "{0}" -f ".md"
$deD2Target = $wsDe.Range("D2").Hyperlinks.Item(1).Address

$wsDe.Range("A2").Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Target, [Type]::Missing, [Type]::Missing, "$newUuid.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $deB2Target, [Type]::Missing, [Type]::Missing, $deB2Display) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deD2Target, [Type]::Missing, [Type]::Missing, "$newUuid.$newHash.de-de.xlf") | Out-Null

$wsDe.Range("E2").Value = "2016-03-13 02:55:52"
